# Headless-Excel COM-interop script implementing:
#  - add defined name "linux"   -> Отчет!$D$6
#  - remove defined names "Тело" and "число_Людей"
#  - add defined name "число1"  -> Отчет!$B$5  (same cell "число_Людей" used to cover)
#  - change selection on sheet "Отчет" to E7
#  - rewrite formulas in B7, C7 and D8 to use the new/renamed defined names

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Отчет")

# --- Defined names -------------------------------------------------------

# "linux" is plain ASCII, so it can be created directly.
$ws.Range("D6").Name = "linux"

# Drop the two obsolete names.
$wb.Names.Item("Тело").Delete()
$wb.Names.Item("число_Людей").Delete()

# "число1" starts with a Cyrillic letter; creating/assigning a name that
# begins with a non-ASCII character directly is not handled correctly by
# this runtime, so create it first under a plain ASCII alias and then
# rename the Name object itself (renaming an existing Name works fine).
$ws.Range("B5").Name = "tmpNameForChislo1"
$wb.Names.Item("tmpNameForChislo1").Name = "число1"

# --- Selection -------------------------------------------------------------

$ws.Range("E7").Select()

# --- Formulas --------------------------------------------------------------

# Row 7: B7 and C7 now reference the (re)defined names instead of summing
# the raw cells directly.
$ws.Range("B7").Formula = "=SUM(число1)"
$ws.Range("C7").Formula = "=SUM(Командировки)"

# Row 8: D8 now sums the new "linux" name instead of being part of the
# C8:E8 shared formula.
$ws.Range("D8").Formula = "=SUM(linux)"

$wb.Save()
